$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 100.07
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 100.07
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2062.87
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 2062.87

# Row 3
$ws.Range("E3").Value = 100.07
$ws.Range("I3").Value = 100.07
$ws.Range("J3").Value = 0.08
$ws.Range("K3").Value = 73.28
$ws.Range("M3").Value = 73.37

# Row 4
$ws.Range("E4").Value = 100.07
$ws.Range("I4").Value = 100.07
$ws.Range("J4").Value = 0.05
$ws.Range("K4").Value = 111.99
$ws.Range("M4").Value = 112.04

# Row 5
$ws.Range("E5").Value = 100.07
$ws.Range("H5").Value = 0.18
$ws.Range("I5").Value = 100.07
$ws.Range("J5").Value = 0.15
$ws.Range("K5").Value = 37.43
$ws.Range("L5").Value = 0.18
$ws.Range("M5").Value = 37.6

# Row 6
$ws.Range("E6").Value = 100.07
$ws.Range("H6").Value = -0.04
$ws.Range("I6").Value = 100.07
$ws.Range("J6").Value = -0.04
$ws.Range("K6").Value = 158.14
$ws.Range("L6").Value = -0.04
$ws.Range("M6").Value = 158.18
